# Daily attendance processing - 2025-12-31 10:59:19
#
# Normalizes the "Recorded By" (column G) lists on the "Session Analysis
# Results" sheet so that any "System"/"system" author(s) sort to the front
# of the comma-separated list (keeping their relative order), while the
# remaining (non-system) authors are moved to the end in reverse order.
# Rows whose list already satisfies this ordering (or that only contain a
# single author) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-RecordedBy($s) {
    if ($null -eq $s -or $s -eq "") {
        return $s
    }

    $parts = $s -split ", "
    if ($parts.Count -le 1) {
        return $s
    }

    $systemItems = @()
    $otherItems = @()
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $systemItems += $p
        } else {
            $otherItems += $p
        }
    }

    # Reverse the non-system items (manual loop - [array]::Reverse() does
    # not mutate reliably in this host).
    $reversedOther = @()
    for ($i = $otherItems.Count - 1; $i -ge 0; $i--) {
        $reversedOther += $otherItems[$i]
    }

    $result = $systemItems + $reversedOther
    return ($result -join ", ")
}

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2
    if ($null -eq $current -or $current -eq "") {
        continue
    }
    if ($current -eq "Recorded By") {
        continue
    }

    $new = Transform-RecordedBy $current
    if ($new -ne $current) {
        $cell.Value2 = $new
    }
}
